$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (D = Price, E = Volume(1h)) per commit diff.
# Leading "'" forces Excel to treat the assigned value as literal text
# (preserves exact formatting: trailing zeros, thousand-dot grouping,
# padding spaces, etc.) without altering the cell NumberFormat.

$ws.Range("D2").Value = "'" + '26.113.08'
$ws.Range("E2").Value = "'" + '  -0.13%  '
$ws.Range("D3").Value = "'" + '1.649.22'
$ws.Range("E3").Value = "'" + '  -1.18%  '
$ws.Range("E4").Value = "'" + '  -0.18%  '
$ws.Range("D5").Value = "'" + '215.93'
$ws.Range("E5").Value = "'" + '  +2.97%  '
$ws.Range("D6").Value = "'" + '0.5222'
$ws.Range("E6").Value = "'" + '  -0.33%  '
$ws.Range("D7").Value = "'" + '1.001'
$ws.Range("E7").Value = "'" + '  -0.15%  '
$ws.Range("D8").Value = "'" + '0.2606'
$ws.Range("E8").Value = "'" + '  -0.68%  '
$ws.Range("D9").Value = "'" + '0.06320'
$ws.Range("E9").Value = "'" + '  -0.38%  '
$ws.Range("D10").Value = "'" + '20.80'
$ws.Range("E10").Value = "'" + '  -1.92%  '
$ws.Range("D11").Value = "'" + '0.07691'
$ws.Range("E11").Value = "'" + '  +2.03%  '
$ws.Range("D12").Value = "'" + '1.647.92'
$ws.Range("E12").Value = "'" + '  -1.38%  '
$ws.Range("D13").Value = "'" + '4.419'
$ws.Range("E13").Value = "'" + '  -0.61%  '
$ws.Range("D14").Value = "'" + '1.868.33'
$ws.Range("E14").Value = "'" + '  -1.58%  '
$ws.Range("D15").Value = "'" + '0.5582'
$ws.Range("E15").Value = "'" + '  +1.35%  '
$ws.Range("D16").Value = "'" + '0.0₅8190'
$ws.Range("D17").Value = "'" + '65.28'
$ws.Range("E17").Value = "'" + '  -1.76%  '
$ws.Range("D18").Value = "'" + '26.104.59'
$ws.Range("E18").Value = "'" + '  -0.22%  '
$ws.Range("D19").Value = "'" + '1.002'
$ws.Range("E19").Value = "'" + '  -0.10%  '
$ws.Range("D20").Value = "'" + '4.739'
$ws.Range("E20").Value = "'" + '  -0.33%  '
$ws.Range("D21").Value = "'" + '188.80'
$ws.Range("E21").Value = "'" + '  +1.11%  '
$ws.Range("D22").Value = "'" + '10.24'
$ws.Range("E22").Value = "'" + '  -0.69%  '
$ws.Range("D23").Value = "'" + '6.213'
$ws.Range("E23").Value = "'" + '  +0.44%  '
$ws.Range("E24").Value = "'" + '  -0.12%  '
$ws.Range("D25").Value = "'" + '145.96'
$ws.Range("E25").Value = "'" + '  -2.36%  '
$ws.Range("D26").Value = "'" + '7.446'
$ws.Range("E26").Value = "'" + '  -0.86%  '
$ws.Range("D27").Value = "'" + '0.1209'
$ws.Range("E27").Value = "'" + '  -3.30%  '
$ws.Range("D28").Value = "'" + '15.88'
$ws.Range("E28").Value = "'" + '  +0.11%  '
$ws.Range("E29").Value = "'" + '  +3.11%  '
$ws.Range("D30").Value = "'" + '0.05921'
$ws.Range("E30").Value = "'" + '  -7.77%  '
$ws.Range("D31").Value = "'" + '1.265'
$ws.Range("E31").Value = "'" + '  -0.85%  '
$ws.Range("D32").Value = "'" + '3.435'
$ws.Range("E32").Value = "'" + '  -2.21%  '
$ws.Range("D33").Value = "'" + '3.409'
$ws.Range("E33").Value = "'" + '  -0.16%  '
$ws.Range("D34").Value = "'" + '1.653'
$ws.Range("E34").Value = "'" + '  +0.52%  '
$ws.Range("D35").Value = "'" + '0.9864'
$ws.Range("E35").Value = "'" + '  -1.99%  '
$ws.Range("D36").Value = "'" + '2.393'
$ws.Range("E36").Value = "'" + '  -0.59%  '
$ws.Range("D37").Value = "'" + '2.754'
$ws.Range("E37").Value = "'" + '  +0.53%  '
$ws.Range("D38").Value = "'" + '0.5678'
$ws.Range("E38").Value = "'" + '  -5.76%  '
$ws.Range("E39").Value = "'" + '  -0.03%  '
$ws.Range("D40").Value = "'" + '0.8602'
$ws.Range("E40").Value = "'" + '  -0.71%  '
$ws.Range("D41").Value = "'" + '5.782'
$ws.Range("E41").Value = "'" + '  -5.91%  '
$ws.Range("E42").Value = "'" + '  -0.24%  '
$ws.Range("D43").Value = "'" + '1.027.08'
$ws.Range("E43").Value = "'" + '  -7.82%  '
$ws.Range("D44").Value = "'" + '100.22'
$ws.Range("E44").Value = "'" + '  -0.13%  '
$ws.Range("D45").Value = "'" + '1.795.83'
$ws.Range("E45").Value = "'" + '  -1.43%  '
$ws.Range("D46").Value = "'" + '0.0₈107'
$ws.Range("E46").Value = "'" + '  -2.74%  '
$ws.Range("D47").Value = "'" + '56.03'
$ws.Range("E47").Value = "'" + '  +0.95%  '
$ws.Range("D48").Value = "'" + '1.002'
$ws.Range("E48").Value = "'" + '  -0.10%  '
$ws.Range("D49").Value = "'" + '8.090'
$ws.Range("E49").Value = "'" + '  +0.16%  '
$ws.Range("D50").Value = "'" + '0.05188'
$ws.Range("E50").Value = "'" + '  -0.85%  '
$ws.Range("D51").Value = "'" + '0.4220'
$ws.Range("E51").Value = "'" + '  -0.45%  '
